$d = $word.ActiveDocument
$s = $d.Styles.Add("Paragraphedeliste", 1)
$pf = $s.ParagraphFormat
$pf.LeftIndent = 36
Write-Output "done"
